$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$refStyle = $ws.Range("B2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "73.294.30"
$ws.Range("D2").Style = $refStyle
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.971.79"
$ws.Range("D3").Style = $refStyle
$ws.Range("E3").Value = "  -1.89%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.79"
$ws.Range("D5").Style = $refStyle
$ws.Range("E5").Value = "  +9.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.45"
$ws.Range("D6").Style = $refStyle
$ws.Range("E6").Value = "  +11.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.681"
$ws.Range("D7").Style = $refStyle
$ws.Range("E7").Value = "  -2.13%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.779"
$ws.Range("D9").Style = $refStyle
$ws.Range("E9").Value = "  +1.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.186"
$ws.Range("D10").Style = $refStyle
$ws.Range("E10").Value = "  +7.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.03"
$ws.Range("D11").Style = $refStyle
$ws.Range("E11").Value = "  +3.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000336"
$ws.Range("D12").Style = $refStyle
$ws.Range("E12").Value = "  +1.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.28"
$ws.Range("D13").Style = $refStyle
$ws.Range("E13").Value = "  +2.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.605.93"
$ws.Range("D14").Style = $refStyle
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.975.43"
$ws.Range("D15").Style = $refStyle
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.19"
$ws.Range("D16").Style = $refStyle
$ws.Range("E16").Value = "  -2.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.24"
$ws.Range("D17").Style = $refStyle
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.46"
$ws.Range("D18").Style = $refStyle
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.200.67"
$ws.Range("D19").Style = $refStyle
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "439.56"
$ws.Range("D21").Style = $refStyle
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.84"
$ws.Range("D22").Style = $refStyle
$ws.Range("E22").Value = "  +9.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "95.65"
$ws.Range("D23").Style = $refStyle
$ws.Range("E23").Value = "  -2.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.38"
$ws.Range("D24").Style = $refStyle
$ws.Range("E24").Value = "  -4.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.17"
$ws.Range("D25").Style = $refStyle
$ws.Range("E25").Value = "  -4.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.16"
$ws.Range("D26").Style = $refStyle
$ws.Range("E26").Value = "  -5.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.07"
$ws.Range("D27").Style = $refStyle
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.96"
$ws.Range("D28").Style = $refStyle
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.49"
$ws.Range("D29").Style = $refStyle
$ws.Range("E29").Value = "  -4.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.10"
$ws.Range("D30").Style = $refStyle
$ws.Range("E30").Value = "  -3.30%  "
$ws.Range("E31").Value = "  -1.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.81"
$ws.Range("D32").Style = $refStyle
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0000106"
$ws.Range("D33").Style = $refStyle
$ws.Range("E33").Value = "  +14.31%  "
$ws.Range("E34").Value = "  -3.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "48.10"
$ws.Range("D35").Style = $refStyle
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "70.46"
$ws.Range("D36").Style = $refStyle
$ws.Range("E36").Value = "  +4.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "648.18"
$ws.Range("D37").Style = $refStyle
$ws.Range("E37").Value = "  -5.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.430"
$ws.Range("D38").Style = $refStyle
$ws.Range("E38").Value = "  -4.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.41"
$ws.Range("D39").Style = $refStyle
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = $refStyle
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.146"
$ws.Range("D41").Style = $refStyle
$ws.Range("E41").Value = "  -3.04%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  -2.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.19"
$ws.Range("D44").Style = $refStyle
$ws.Range("E44").Value = "  -5.79%  "
$ws.Range("E45").Value = "  -4.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.12"
$ws.Range("D46").Style = $refStyle
$ws.Range("E46").Value = "  +32.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.149"
$ws.Range("D47").Style = $refStyle
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000299"
$ws.Range("D48").Style = $refStyle
$ws.Range("E48").Value = "  +5.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.43"
$ws.Range("D49").Style = $refStyle
$ws.Range("E49").Value = "  +2.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.57"
$ws.Range("D50").Style = $refStyle
$ws.Range("E50").Value = "  -5.23%  "
$ws.Range("E51").Value = "  -5.08%  "
